# The commit appends " (Changed main)" to the end of the first paragraph's
# sentence "This is a Microsoft word document." so the paragraph reads:
#   "This is a Microsoft word document. (Changed main)"
#
# Use Find/Replace (wdReplaceAll-ish single Execute) on the whole document
# story to retarget the existing sentence onto the expanded sentence.

$d = $word.ActiveDocument

$old = "This is a Microsoft word document."
$new = "This is a Microsoft word document. (Changed main)"

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
